# "Add files via upload" — updates the İş Takip Listesi sheet so that rows
# 123-176 in the "BİRİM" column (D) are re-labelled from "GÜNCELLEME" to the
# new "GÜNCELLEME(2026)" category (mirroring the newly-added
# "Güncelleme(2026)" worksheet), and restores the on-screen selection over
# that same block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("İş Takip Listesi")
$ws.Activate()

# Re-point D123:D176 at the "GÜNCELLEME(2026)" label (a brand-new shared
# string; everything else in column D above/below this block keeps the
# original "GÜNCELLEME" value).
$ws.Range("D123:D176").Value = "GÜNCELLEME(2026)"

# Match the saved selection/active cell from the edited workbook.
$ws.Range("D123:D176").Select() | Out-Null
